# Sync attendance_reports: normalize "Recorded By" ordering in column G.
#
# For a handful of rows the first two comma-separated names/emails in the
# "Recorded By" cell (column G) were swapped, e.g.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, admin@admin.com"             -> "admin@admin.com, System"
#   "System, system, backup@backdoor.com" -> "system, System, backup@backdoor.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"
# Other "Recorded By" values (e.g. "System, backup@backdoor.com",
# single-name values, counts like "0/52", etc.) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Exact before -> after replacements observed for the "Recorded By" column.
$replacements = @{
    "System, system, backup@backdoor.com" = "system, System, backup@backdoor.com";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($null -ne $current -and $replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
